# Applies the scheduled-runner profit-recalculation update across the
# Behemoth_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 4434.375
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 4434.375
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 4434.375
$ws.Cells.Item(40, 13).Value = ""
$ws.Cells.Item(40, 14).Value = -4784.375

$ws.Cells.Item(64, 8).Value = 5079.6665
$ws.Cells.Item(64, 10).Value = 4981.8184
$ws.Cells.Item(64, 12).Value = 4981.8184
$ws.Cells.Item(64, 14).Value = -5477.8184

$ws.Cells.Item(67, 8).Value = 5079.6665
$ws.Cells.Item(67, 10).Value = 4981.8184
$ws.Cells.Item(67, 12).Value = 4981.8184
$ws.Cells.Item(67, 14).Value = -6697.8184

$ws.Cells.Item(80, 8).Value = 824.75
$ws.Cells.Item(80, 9).Value = 449.5
$ws.Cells.Item(80, 10).Value = 1200
$ws.Cells.Item(80, 11).Value = 1348.5
$ws.Cells.Item(80, 12).Value = 3600
$ws.Cells.Item(80, 13).Value = -350.5
$ws.Cells.Item(80, 14).Value = -5596

$ws.Cells.Item(83, 8).Value = 824.75
$ws.Cells.Item(83, 9).Value = 449.5
$ws.Cells.Item(83, 10).Value = 1200
$ws.Cells.Item(83, 11).Value = 4045.5
$ws.Cells.Item(83, 12).Value = 10800
$ws.Cells.Item(83, 13).Value = 946.5
$ws.Cells.Item(83, 14).Value = -20784

$ws.Cells.Item(86, 8).Value = 6689.5654
$ws.Cells.Item(86, 9).Value = 6937
$ws.Cells.Item(86, 11).Value = 6937
$ws.Cells.Item(86, 13).Value = -5814

$ws.Cells.Item(89, 8).Value = 6689.5654
$ws.Cells.Item(89, 9).Value = 6937
$ws.Cells.Item(89, 11).Value = 34685
$ws.Cells.Item(89, 13).Value = -29069

$ws.Cells.Item(106, 8).Value = 2941.8572
$ws.Cells.Item(106, 9).Value = 2498.8333
$ws.Cells.Item(106, 11).Value = 2498.8333
$ws.Cells.Item(106, 13).Value = -1867.8333

$ws.Cells.Item(137, 8).Value = 5226.6924
$ws.Cells.Item(137, 9).Value = 1777.4445
$ws.Cells.Item(137, 10).Value = 12987.5
$ws.Cells.Item(137, 11).Value = 5332.333500000001
$ws.Cells.Item(137, 12).Value = 38962.5
$ws.Cells.Item(137, 13).Value = -2782.333500000001
$ws.Cells.Item(137, 14).Value = -44062.5

$ws.Cells.Item(138, 8).Value = 2864.577
$ws.Cells.Item(138, 9).Value = 1022.6667
$ws.Cells.Item(138, 10).Value = 3104.8262
$ws.Cells.Item(138, 11).Value = 3068.0001
$ws.Cells.Item(138, 12).Value = 9314.4786
$ws.Cells.Item(138, 13).Value = 2071.9999
$ws.Cells.Item(138, 14).Value = -19594.4786

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 2999
$ws.Cells.Item(63, 9).Value = 2999
$ws.Cells.Item(63, 11).Value = 2999
$ws.Cells.Item(63, 13).Value = -2313

$ws.Cells.Item(66, 8).Value = 2999
$ws.Cells.Item(66, 9).Value = 2999
$ws.Cells.Item(66, 11).Value = 14995
$ws.Cells.Item(66, 13).Value = -11563

$ws.Cells.Item(122, 8).Value = 1235
$ws.Cells.Item(122, 9).Value = 1141.3636
$ws.Cells.Item(122, 10).Value = 1750
$ws.Cells.Item(122, 11).Value = 3424.0908
$ws.Cells.Item(122, 12).Value = 5250
$ws.Cells.Item(122, 13).Value = -974.0907999999999
$ws.Cells.Item(122, 14).Value = -10150

$ws.Cells.Item(132, 8).Value = 3340
$ws.Cells.Item(132, 9).Value = 2127.7778
$ws.Cells.Item(132, 11).Value = 6383.3334
$ws.Cells.Item(132, 13).Value = -3853.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 1365.8334
$ws.Cells.Item(22, 9).Value = 1232.6666
$ws.Cells.Item(22, 11).Value = 1232.6666
$ws.Cells.Item(22, 13).Value = -1059.6666

$ws.Cells.Item(35, 8).Value = 67
$ws.Cells.Item(35, 9).Value = 67
$ws.Cells.Item(35, 11).Value = 67
$ws.Cells.Item(35, 13).Value = 243

$ws.Cells.Item(94, 8).Value = 1038.9062
$ws.Cells.Item(94, 9).Value = 1054.6774
$ws.Cells.Item(94, 10).Value = 550
$ws.Cells.Item(94, 11).Value = 1054.6774
$ws.Cells.Item(94, 12).Value = 550
$ws.Cells.Item(94, 13).Value = -603.6774
$ws.Cells.Item(94, 14).Value = -1452

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 2723.1538
$ws.Cells.Item(7, 10).Value = 6883
$ws.Cells.Item(7, 12).Value = 6883
$ws.Cells.Item(7, 14).Value = -7109

$ws.Cells.Item(31, 8).Value = 708970.8
$ws.Cells.Item(31, 9).Value = 1434.0555
$ws.Cells.Item(31, 11).Value = 1434.0555
$ws.Cells.Item(31, 13).Value = -1139.0555

$ws.Cells.Item(34, 8).Value = 708970.8
$ws.Cells.Item(34, 9).Value = 1434.0555
$ws.Cells.Item(34, 11).Value = 1434.0555
$ws.Cells.Item(34, 13).Value = -1232.0555

$ws.Cells.Item(132, 8).Value = 2110.25
$ws.Cells.Item(132, 9).Value = 2158.5217
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 6475.5651
$ws.Cells.Item(132, 12).Value = 3000
$ws.Cells.Item(132, 13).Value = -3945.5651
$ws.Cells.Item(132, 14).Value = -8060

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 1000
$ws.Cells.Item(23, 10).Value = 1000
$ws.Cells.Item(23, 12).Value = 3000
$ws.Cells.Item(23, 14).Value = -3470

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 3174.375
$ws.Cells.Item(3, 10).Value = 799.1667
$ws.Cells.Item(3, 12).Value = 799.1667
$ws.Cells.Item(3, 14).Value = -1031.1667

$ws.Cells.Item(70, 8).Value = 4740.7144
$ws.Cells.Item(70, 9).Value = 4740.7144
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 4740.7144
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = -4470.7144
$ws.Cells.Item(70, 14).Value = ""

$ws.Cells.Item(73, 8).Value = 4740.7144
$ws.Cells.Item(73, 9).Value = 4740.7144
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 4740.7144
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = -3804.7144
$ws.Cells.Item(73, 14).Value = ""

$ws.Cells.Item(80, 8).Value = 440
$ws.Cells.Item(80, 9).Value = 440
$ws.Cells.Item(80, 11).Value = 440
$ws.Cells.Item(80, 13).Value = 558

$ws.Cells.Item(83, 8).Value = 440
$ws.Cells.Item(83, 9).Value = 440
$ws.Cells.Item(83, 11).Value = 2200
$ws.Cells.Item(83, 13).Value = 2792

$ws.Cells.Item(97, 8).Value = 1457
$ws.Cells.Item(97, 9).Value = 1074.375
$ws.Cells.Item(97, 11).Value = 1074.375
$ws.Cells.Item(97, 13).Value = -578.375

$ws.Cells.Item(122, 8).Value = 1503.24
$ws.Cells.Item(122, 9).Value = 1463.409
$ws.Cells.Item(122, 11).Value = 4390.227000000001
$ws.Cells.Item(122, 13).Value = -1940.227000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3476.1428
$ws.Cells.Item(46, 9).Value = 3450
$ws.Cells.Item(46, 10).Value = 3559.8
$ws.Cells.Item(46, 11).Value = 3450
$ws.Cells.Item(46, 12).Value = 3559.8
$ws.Cells.Item(46, 13).Value = -3262
$ws.Cells.Item(46, 14).Value = -3935.8

$ws.Cells.Item(136, 8).Value = 76148.625
$ws.Cells.Item(136, 9).Value = 5037.4
$ws.Cells.Item(136, 11).Value = 15112.2
$ws.Cells.Item(136, 13).Value = -12562.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 150500
$ws.Cells.Item(5, 9).Value = 1000
$ws.Cells.Item(5, 11).Value = 1000
$ws.Cells.Item(5, 13).Value = -888

$ws.Cells.Item(51, 8).Value = 33534
$ws.Cells.Item(51, 9).Value = 32070
$ws.Cells.Item(51, 11).Value = 32070
$ws.Cells.Item(51, 13).Value = -31560

$ws.Cells.Item(55, 8).Value = 1498.5
$ws.Cells.Item(55, 9).Value = 1498.5
$ws.Cells.Item(55, 11).Value = 1498.5
$ws.Cells.Item(55, 13).Value = -1221.5

$ws.Cells.Item(75, 8).Value = 18814244
$ws.Cells.Item(75, 10).Value = 18814244
$ws.Cells.Item(75, 12).Value = 18814244
$ws.Cells.Item(75, 14).Value = -18816116

$ws.Cells.Item(78, 8).Value = 18814244
$ws.Cells.Item(78, 10).Value = 18814244
$ws.Cells.Item(78, 12).Value = 56442732
$ws.Cells.Item(78, 14).Value = -56452092

$ws.Cells.Item(95, 8).Value = 66072.75
$ws.Cells.Item(95, 10).Value = 66072.75
$ws.Cells.Item(95, 12).Value = 66072.75
$ws.Cells.Item(95, 14).Value = -71564.75

$ws.Cells.Item(107, 8).Value = 17858570
$ws.Cells.Item(107, 9).Value = 22728754
$ws.Cells.Item(107, 11).Value = 68186262
$ws.Cells.Item(107, 13).Value = -68184342

Write-Host "Applied 185 cell updates across 38 rows."
